$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Structural change: insert two new columns -----------------------------
# New column C: "MinCount" (=SUM of MinCount_ADULTMF + MinCount_CALFMF)
$ws.Range("C1").EntireColumn.Insert()
# New column O (was N before the first insert): "Estimate" (=SUM of Estimate_ADULTMF + Estimate_CALFMF)
$ws.Range("O1").EntireColumn.Insert()

# --- New header labels -------------------------------------------------------
# Write O1 ("Estimate") before C1 ("MinCount") so new shared strings land in
# the same order as the source workbook.
$ws.Range("O1").Value = "Estimate"
$ws.Range("C1").Value = "MinCount"

# --- Row 2 (2001): sex ratio column now populated ---------------------------
$ws.Range("H2").Value = 0.64

# --- Row 3 (2002) ------------------------------------------------------------
$ws.Range("B3").Formula = "=G3+F3"
$ws.Range("C3").Formula = "=SUM(D3:E3)"
$ws.Range("J3").Value = 0.22
$ws.Range("O3").Formula = "=SUM(M3,K3)"

# --- Row 4 (2003) -------------------------------------------------------------
$ws.Range("H4").Value = 0.64

# --- Row 5 (2004) -------------------------------------------------------------
$ws.Range("H5").Value = 0.64

# --- Row 6 (2005) -------------------------------------------------------------
$ws.Range("H6").Value = 0.64

# --- Row 7 (2006) -------------------------------------------------------------
$ws.Range("H7").Value = 0.64

# --- Row 8 (2007) -------------------------------------------------------------
$ws.Range("H8").Value = 0.64

# --- Row 9 (2008) -------------------------------------------------------------
$ws.Range("C9").Formula = "=SUM(D9:E9)"
$ws.Range("J9").Value = 0.13
$ws.Range("O9").Formula = "=SUM(M9,K9)"

# --- Row 10 (2009) ------------------------------------------------------------
$ws.Range("H10").Value = 0.64

# --- Row 11 (2010) ------------------------------------------------------------
$ws.Range("H11").Value = 0.64

# --- Row 12 (2011) ------------------------------------------------------------
$ws.Range("H12").Value = 0.64

# --- Row 13 (2012) ------------------------------------------------------------
$ws.Range("H13").Value = 0.64

# --- Row 14 (2013) ------------------------------------------------------------
$ws.Range("C14").Formula = "=SUM(D14:E14)"
$ws.Range("I14").Value = 0.85
$ws.Range("J14").Value = 0.1
$ws.Range("O14").Formula = "=SUM(M14,K14)"

# --- Row 15 (2014) ------------------------------------------------------------
$ws.Range("H15").Value = 0.64

# --- Row 16 (2015) ------------------------------------------------------------
$ws.Range("H16").Value = 0.64

# --- Row 17 (2016) -- data update: juveniles corrected from 6 to 5 -----------
$ws.Range("C17").Formula = "=SUM(D17:E17)"
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 5
$ws.Range("J17").Value = 0.13
$ws.Range("O17").Formula = "=SUM(M17,K17)"
$ws.Range("Q17").Value = "Slight difference from the report due to removal of one group of 2 caribou (with one collar) that was outside the herd range . Data Update: 2016 report records 'juveniles = 6', updated to 'juveniles = 5', one group of two caribou (one collared, car 191) determined to be outside herd boundary."

# --- Row 18 (2017) ------------------------------------------------------------
$ws.Range("H18").Value = 0.64

# --- Row 19 (2018) -- MinCount updated with KMB gov numbers ------------------
$ws.Range("C19").Formula = "=SUM(D19:E19)"
$ws.Range("D19").Value = 54
$ws.Range("E19").Value = 13
$ws.Range("O19").Formula = "=SUM(M19,K19)"
$ws.Range("Q19").Value = "#'s from KMB updated gov #'s"

# --- Row 20 (2019) -- MinCount updated with KMB gov numbers ------------------
$ws.Range("C20").Formula = "=SUM(D20:E20)"
$ws.Range("O20").Formula = "=SUM(M20,K20)"
$ws.Range("Q20").Value = "#'s from KMB updated gov #'s"

# --- Row 21 (2020) ------------------------------------------------------------
$ws.Range("H21").Value = 0.64

# --- Row 22 (2021) ------------------------------------------------------------
$ws.Range("B22").Formula = "=G22+F22"
$ws.Range("C22").Formula = "=SUM(D22:E22)"
$ws.Range("O22").Formula = "=SUM(M22,K22)"

# --- Selection / view bookkeeping (matches the author's final click) ---------
$ws.Range("E27").Select()
